# Applies the "Import 6 months BB credit card history" data-sync update
# to Dashboard_2026.xlsx:
#   - Refresh sync timestamps (Dashboard!A2, Dados!B3)
#   - Dashboard: update Gastos Variaveis summary + per-category rows
#   - Mensal: zero out the old single-month Budget column, add a new
#     "Jun real" column (M) with per-category totals
#   - Categorias: turn the header placeholders for Gasto Real/Disponivel/%
#     into live numbers, refresh per-category rows
#   - Dados: refresh the synced gasto_jan snapshot column

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Dashboard
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A2").Value = "Atualizado: 30/12/2025 18:33"

# Gastos Variaveis summary row
$dash.Range("B7").Value = 21000
$dash.Range("C7").Value = 29701.39
# D7 holds literal text ("41%"), not a percentage-formatted number - force
# Text format before the write so Excel doesn't auto-convert it to 0.41,
# then restore the original General format (copied from the untouched D6
# cell) so the cell's style index is unchanged.
$dash.Range("D7").NumberFormat = "@"
$dash.Range("D7").Value = "41%"
$dash.Range("D6").Copy()
$dash.Range("D7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Obra row - only "Real" changes
$dash.Range("C9").Value = 15798.61

# GASTOS POR CATEGORIA table
$dash.Range("C14").Value = 5048.87
$dash.Range("D14").Value = 126

$dash.Range("B15").Value = 3500
$dash.Range("C15").Value = 8210.299999999999
$dash.Range("D15").Value = 234

$dash.Range("B16").Value = 2000
$dash.Range("C16").Value = 9500
$dash.Range("D16").Value = 475

$dash.Range("B17").Value = 4200
$dash.Range("C17").Value = 222.91
$dash.Range("D17").Value = 5

$dash.Range("B18").Value = 3800
$dash.Range("C18").Value = 567.4
$dash.Range("D18").Value = 14

$dash.Range("B19").Value = 1300
$dash.Range("C19").Value = 5493.15
$dash.Range("D19").Value = 422

$dash.Range("B20").Value = 1500
$dash.Range("C20").Value = 618.86
$dash.Range("D20").Value = 41

$dash.Range("B21").Value = 400
$dash.Range("C21").Value = 39.9
$dash.Range("D21").Value = 9

$dash.Range("B22").Value = 300

# ---------------------------------------------------------------
# Mensal (6-month variable-spend view)
# ---------------------------------------------------------------
$mensal = $wb.Worksheets.Item("Mensal")

$mensal.Range("B4").Value = 0
$mensal.Range("M4").Value = 5048.87

$mensal.Range("B5").Value = 0
$mensal.Range("M5").Value = 8210.299999999999

$mensal.Range("B6").Value = 0
$mensal.Range("M6").Value = 9500

$mensal.Range("B7").Value = 0
$mensal.Range("M7").Value = 222.91

$mensal.Range("B8").Value = 0
$mensal.Range("M8").Value = 567.4

$mensal.Range("B9").Value = 0
$mensal.Range("M9").Value = 5493.15

$mensal.Range("B10").Value = 0
$mensal.Range("M10").Value = 618.86

$mensal.Range("B11").Value = 0
$mensal.Range("M11").Value = 39.9

$mensal.Range("B12").Value = 0
$mensal.Range("M12").Value = 0

# ---------------------------------------------------------------
# Categorias
# ---------------------------------------------------------------
$cat = $wb.Worksheets.Item("Categorias")

$cat.Range("C4").Value = 5048.87
$cat.Range("D4").Value = -1048.87
$cat.Range("E4").Value = 1.2622175

$cat.Range("C5").Value = 8210.299999999999
$cat.Range("D5").Value = -4710.299999999999
$cat.Range("E5").Value = 2.3458

$cat.Range("C6").Value = 9500
$cat.Range("D6").Value = -7500
$cat.Range("E6").Value = 4.75

$cat.Range("C7").Value = 222.91
$cat.Range("D7").Value = 3977.09
$cat.Range("E7").Value = 0.05307380952380952

$cat.Range("C8").Value = 567.4
$cat.Range("D8").Value = 3232.6
$cat.Range("E8").Value = 0.1493157894736842

$cat.Range("C9").Value = 5493.15
$cat.Range("D9").Value = -4193.15
$cat.Range("E9").Value = 4.225499999999999

$cat.Range("C10").Value = 618.86
$cat.Range("D10").Value = 881.14
$cat.Range("E10").Value = 0.4125733333333333

$cat.Range("C11").Value = 39.9
$cat.Range("D11").Value = 360.1
$cat.Range("E11").Value = 0.09974999999999999

$cat.Range("D12").Value = 300

# ---------------------------------------------------------------
# Dados (synced snapshot)
# ---------------------------------------------------------------
$dados = $wb.Worksheets.Item("Dados")

$dados.Range("B3").Value = "2025-12-30T18:33:45.564073"

$dados.Range("D8").Value = 5048.87
$dados.Range("D9").Value = 8210.299999999999
$dados.Range("D10").Value = 9500
$dados.Range("D11").Value = 222.91
$dados.Range("D12").Value = 567.4
$dados.Range("D13").Value = 5493.15
$dados.Range("D14").Value = 618.86
$dados.Range("D15").Value = 39.9

Write-Output "edit.ps1 applied"
